$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.728.31"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "2.516.85"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'534.73"
$ws.Range("E5").Value = "  +5.69%  "
$ws.Range("D6").Value = "'134.11"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("D9").Value = "2.515.02"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").Value = "'0.0994"
$ws.Range("E10").Value = "  +4.31%  "
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "2.952.63"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").Value = "58.710.81"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "'22.33"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("E17").Value = "  +3.03%  "
$ws.Range("D18").Value = "2.502.24"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").Value = "'320.75"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  +9.20%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'65.40"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "'0.995"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'7.49"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("D29").Value = "0.0₃0762"
$ws.Range("E29").Value = "  +5.46%  "
$ws.Range("D30").Value = "'172.21"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("E31").Value = "  +5.48%  "
$ws.Range("E32").Value = "  +4.89%  "
$ws.Range("D33").Value = "'6.28"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'18.14"
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").Value = "'3.95"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("E39").Value = "  +4.24%  "
$ws.Range("D40").Value = "'36.68"
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").Value = "'0.800"
$ws.Range("E41").Value = "  +5.45%  "
$ws.Range("D42").Value = "'5.18"
$ws.Range("E42").Value = "  +6.68%  "
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("D44").Value = "'276.52"
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("E45").Value = "  +10.01%  "
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").Value = "'0.0935"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("D48").Value = "'0.0509"
$ws.Range("E48").Value = "  +5.35%  "
$ws.Range("E49").Value = "  +4.95%  "
$ws.Range("D50").Value = "'17.02"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").Value = "1.753.68"
$ws.Range("E51").Value = "  +2.97%  "
